$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained a new weekly data point (2022-10-07, serial 44841) for
# "Apio" / "Americana (o)" at Vega Central Mapocho de Santiago, reported in
# both "Primera" and "Segunda" quality grades. These two rows are inserted
# right after the existing row 259, pushing all following rows down by two
# (rows 260..336 become 262..338).

# Insert two fresh rows at position 260 (each Insert() shifts everything
# at/after row 260 down by one, mirroring Excel's native "Insert Rows").
$ws.Rows.Item(260).Insert()
$ws.Rows.Item(260).Insert()

# New row 260: Apio / Americana (o) / Primera
$ws.Cells.Item(260,1).Value = 9
$ws.Cells.Item(260,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(260,3).Value = "Metropolitana"
$ws.Cells.Item(260,4).Value = 44841
$ws.Cells.Item(260,5).Value = 13
$ws.Cells.Item(260,6).Value = 100112017
$ws.Cells.Item(260,7).Value = "Apio"
$ws.Cells.Item(260,8).Value = "Americana (o)"
$ws.Cells.Item(260,9).Value = "Primera"
$ws.Cells.Item(260,10).Value = 90
$ws.Cells.Item(260,11).Value = 9000
$ws.Cells.Item(260,12).Value = 9000
$ws.Cells.Item(260,13).Value = 9000
$ws.Cells.Item(260,14).Value = "$/docena de matas"
$ws.Cells.Item(260,15).Value = "Región de Coquimbo"
$ws.Cells.Item(260,16).Value = 1500
$ws.Cells.Item(260,17).Value = 6
$ws.Cells.Item(260,18).Value = "Hortaliza"

# New row 261: Apio / Americana (o) / Segunda
$ws.Cells.Item(261,1).Value = 9
$ws.Cells.Item(261,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(261,3).Value = "Metropolitana"
$ws.Cells.Item(261,4).Value = 44841
$ws.Cells.Item(261,5).Value = 13
$ws.Cells.Item(261,6).Value = 100112017
$ws.Cells.Item(261,7).Value = "Apio"
$ws.Cells.Item(261,8).Value = "Americana (o)"
$ws.Cells.Item(261,9).Value = "Segunda"
$ws.Cells.Item(261,10).Value = 36
$ws.Cells.Item(261,11).Value = 7000
$ws.Cells.Item(261,12).Value = 7000
$ws.Cells.Item(261,13).Value = 7000
$ws.Cells.Item(261,14).Value = "$/docena de matas"
$ws.Cells.Item(261,15).Value = "Región de Coquimbo"
$ws.Cells.Item(261,16).Value = 1167
$ws.Cells.Item(261,17).Value = 6
$ws.Cells.Item(261,18).Value = "Hortaliza"
